$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Booking Number"
$ws.Range("B1").Value = "Conf Number "
$ws.Range("C1").Value = "Guest Name"
$ws.Range("D1").Value = "Price"
